$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Bad Chart day!" cell (G15) to include a link to the slides.
$ws.Range("G15").Value = "- [Bad Chart day!](../slides/14-graphics.qmd)"

# 2. Add the Homework 13 assignment text into H16, with the "(" and ")"
#    around "Submitted via Github classroom" rendered in bold - matching
#    the style used by the other "(Submitted via Github classroom)" cells.
$hw13 = "- [ ] [Homework 13](../homework/13-practice-final.qmd) (Submitted via Github classroom)"
$ws.Range("H16").Value = $hw13
$ws.Range("H16").Characters(56,1).Font.Bold = $true
$ws.Range("H16").Characters(87,1).Font.Bold = $true

# 3. Move the active selection to H17 (matches the saved selection state).
[void]$ws.Range("H17").Select()
